$d = $word.ActiveDocument

# The target paragraph is the last paragraph in the document body (just
# before the section properties), which currently contains nothing but
# the hidden "_GoBack" bookmark left over from the last edit location.
$lastIndex = $d.Paragraphs.Count
$para = $d.Paragraphs.Item($lastIndex)

# Turn it into a numbered/bulleted list item at the same level/list as
# the earlier bullets in this document (numId=1, ilvl=0) by reusing that
# paragraph's list template and telling Word to continue the existing
# list rather than starting a new one.
$firstListPara = $d.Paragraphs.Item(2)
$listTemplate = $firstListPara.Range.ListFormat.ListTemplate
$para.Range.ListFormat.ApplyListTemplateWithLevel($listTemplate, $true, 0, $false, $false)

# Insert the trailing sentence first, via InsertAfter on the untouched
# paragraph range, so it lands after the existing bookmark.
$para.Range.InsertAfter(" is held in the prototypes?")

# Now insert the leading sentence at the very start of the paragraph
# (still before the bookmark and before the run we just added).
$para2 = $d.Paragraphs.Item($lastIndex)
$startRng = $d.Range($para2.Range.Start, $para2.Range.Start)
$startRng.InsertBefore("Is it fair to say the constructor holds the data structures?  Behavior")
